$d = $word.ActiveDocument

$replacements = @(
    @("633×4=", "125×5="),
    @("476×5=", "760×8="),
    @("640×8=", "788×8="),
    @("887×7=", "888×7="),
    @("266×9=", "102×7="),
    @("846×7=", "874×6="),
    @("305×9=", "887×3="),
    @("781×3=", "770×2="),
    @("666×6=", "925×6="),
    @("490×9=", "300×7="),
    @("880×3=", "759×4="),
    @("552×7=", "718×5="),
    @("193×2=", "467×2="),
    @("499×3=", "711×2="),
    @("368×3=", "738×9="),
    @("346×6=", "127×9="),
    @("139×4=", "317×7="),
    @("233×7=", "779×3="),
    @("690×8=", "265×9="),
    @("816×7=", "333×5="),
    @("788×5=", "985×7="),
    @("972×7=", "359×5="),
    @("321×6=", "714×6="),
    @("364×8=", "318×8="),
    @("384×3=", "404×7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
